$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "applyFill" style that was left on D7:D12 - reset them to
# the default (unstyled) look so the style is no longer referenced.
$ws.Range("D7:D12").Style = "Normal"

# Add new row 20 - test case sc16
$ws.Range("A20").Value = "sc16"
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = "A sublimit and a restriction on one of two policies"
$ws.Range("F20").Value = "complete"
$ws.Range("G20").Value = "yes"
$ws.Range("H20").Value = "done"

# Update selection / view to match the saved workbook state
$ws.Range("F20").Select()
